$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.808.62'
$ws.Range("E2").Value = '  -2.96%  '

$ws.Range("D3").Value = '3.488.22'
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.99%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.39%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").Value = '3.485.12'
$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.594'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.131'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.84'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.32%  '

$ws.Range("E12").Value = '  -4.51%  '

$ws.Range("D13").Value = '4.091.48'
$ws.Range("E13").Value = '  +0.23%  '

$ws.Range("E14").Value = '  +0.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.87'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.99%  '

$ws.Range("D16").Value = '65.916.00'
$ws.Range("E16").Value = '  -2.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000171'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.76%  '

$ws.Range("D18").Value = '3.490.19'
$ws.Range("E18").Value = '  +0.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '367.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.76'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.56%  '

$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '

$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.539'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("E26").Value = '  -0.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.177'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.59%  '

$ws.Range("E29").Value = '  +0.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '24.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.77'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.13%  '

$ws.Range("E32").Value = '  -3.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.29'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.05'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.55'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '29.58'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +12.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '159.11'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.73%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.887'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.78'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.07%  '

$ws.Range("D41").Value = '2.794.61'
$ws.Range("E41").Value = '  +1.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.53'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -12.34%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.53%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0687'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.87'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.19%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.97%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0288'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '306.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.822'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.98%  '

$ws.Range("E51").Value = '  -4.01%  '
